$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.4723746666666667
$ws.Range("H2").Value2 = 1.417124
$ws.Range("I2").Value2 = 0.02676815992729067
$ws.Range("J2").Value2 = 0.02676815992729067
$ws.Range("M2").Value2 = 159.4836373333333
$ws.Range("N2").Value2 = 478.450912
$ws.Range("O2").Value2 = 0.2983285084902258
$ws.Range("P2").Value2 = 0.2983285084902258
$ws.Range("Q2").Value2 = 75.3360300241209
$ws.Range("R2").Value2 = 678.0242702170881
$ws.Range("S2").Value2 = 0.007985705226136456
$ws.Range("T2").Value2 = 0.007985705226136456

$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.4723746666666667
$ws.Range("H3").Value2 = 1.417124
$ws.Range("I3").Value2 = 0.02676815992729067
$ws.Range("J3").Value2 = 0.02676815992729067
$ws.Range("O3").Value2 = 0.3227862111630279
$ws.Range("P3").Value2 = 0.3227862111630279
$ws.Range("Q3").Value2 = 81.51226250087599
$ws.Range("R3").Value2 = 733.610362507884
$ws.Range("S3").Value2 = 0.008640392922736148
$ws.Range("T3").Value2 = 0.008640392922736148

$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.4723746666666667
$ws.Range("H4").Value2 = 1.417124
$ws.Range("I4").Value2 = 0.02676815992729067
$ws.Range("J4").Value2 = 0.02676815992729067
$ws.Range("M4").Value2 = 74.38770566666666
$ws.Range("N4").Value2 = 223.163117
$ws.Range("O4").Value2 = 0.1391489036280481
$ws.Range("P4").Value2 = 0.1391489036280482
$ws.Range("Q4").Value2 = 35.13886766838978
$ws.Range("R4").Value2 = 316.249809015508
$ws.Range("S4").Value2 = 0.003724760106022749
$ws.Range("T4").Value2 = 0.00372476010602275

$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.4723746666666667
$ws.Range("H5").Value2 = 1.417124
$ws.Range("I5").Value2 = 0.02676815992729067
$ws.Range("J5").Value2 = 0.02676815992729067
$ws.Range("M5").Value2 = 58.41461433333333
$ws.Range("N5").Value2 = 175.243843
$ws.Range("O5").Value2 = 0.1092697975759847
$ws.Range("P5").Value2 = 0.1092697975759848
$ws.Range("Q5").Value2 = 27.59358397417022
$ws.Range("R5").Value2 = 248.342255767532
$ws.Range("S5").Value2 = 0.002924951416736638
$ws.Range("T5").Value2 = 0.002924951416736638

$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.4723746666666667
$ws.Range("H6").Value2 = 1.417124
$ws.Range("I6").Value2 = 0.02676815992729067
$ws.Range("J6").Value2 = 0.02676815992729067
$ws.Range("M6").Value2 = 69.746216
$ws.Range("N6").Value2 = 209.238648
$ws.Range("O6").Value2 = 0.1304665791427133
$ws.Range("P6").Value2 = 0.1304665791427133
$ws.Range("Q6").Value2 = 32.94634553426133
$ws.Range("R6").Value2 = 296.517109808352
$ws.Range("S6").Value2 = 0.003492350255658675
$ws.Range("T6").Value2 = 0.003492350255658676

$ws.Range("I7").Value2 = 0.9656838605972748
$ws.Range("J7").Value2 = 0.9656838605972748
$ws.Range("M7").Value2 = 159.4836373333333
$ws.Range("N7").Value2 = 478.450912
$ws.Range("O7").Value2 = 0.2983285084902258
$ws.Range("P7").Value2 = 0.2983285084902258
$ws.Range("Q7").Value2 = 2717.810582175819
$ws.Range("R7").Value2 = 24460.29523958237
$ws.Range("S7").Value2 = 0.2880910258050681
$ws.Range("T7").Value2 = 0.2880910258050681

$ws.Range("I8").Value2 = 0.9656838605972748
$ws.Range("J8").Value2 = 0.9656838605972748
$ws.Range("O8").Value2 = 0.3227862111630279
$ws.Range("P8").Value2 = 0.3227862111630279
$ws.Range("S8").Value2 = 0.31170943454348
$ws.Range("T8").Value2 = 0.31170943454348

$ws.Range("I9").Value2 = 0.9656838605972748
$ws.Range("J9").Value2 = 0.9656838605972748
$ws.Range("M9").Value2 = 74.38770566666666
$ws.Range("N9").Value2 = 223.163117
$ws.Range("O9").Value2 = 0.1391489036280481
$ws.Range("P9").Value2 = 0.1391489036280482
$ws.Range("Q9").Value2 = 1267.66417561754
$ws.Range("R9").Value2 = 11408.97758055786
$ws.Range("S9").Value2 = 0.1343738504534117
$ws.Range("T9").Value2 = 0.1343738504534117

$ws.Range("I10").Value2 = 0.9656838605972748
$ws.Range("J10").Value2 = 0.9656838605972748
$ws.Range("M10").Value2 = 58.41461433333333
$ws.Range("N10").Value2 = 175.243843
$ws.Range("O10").Value2 = 0.1092697975759847
$ws.Range("P10").Value2 = 0.1092697975759848
$ws.Range("Q10").Value2 = 995.4617266286195
$ws.Range("R10").Value2 = 8959.155539657577
$ws.Range("S10").Value2 = 0.1055200799698597
$ws.Range("T10").Value2 = 0.1055200799698597

$ws.Range("I11").Value2 = 0.9656838605972748
$ws.Range("J11").Value2 = 0.9656838605972748
$ws.Range("M11").Value2 = 69.746216
$ws.Range("N11").Value2 = 209.238648
$ws.Range("O11").Value2 = 0.1304665791427133
$ws.Range("P11").Value2 = 0.1304665791427133
$ws.Range("Q11").Value2 = 1188.567097421608
$ws.Range("R11").Value2 = 10697.10387679447
$ws.Range("S11").Value2 = 0.1259894698254553
$ws.Range("T11").Value2 = 0.1259894698254553

$ws.Range("G12").Value2 = 0.1331983333333333
$ws.Range("H12").Value2 = 0.399595
$ws.Range("I12").Value2 = 0.007547979475434553
$ws.Range("J12").Value2 = 0.007547979475434553
$ws.Range("M12").Value2 = 159.4836373333333
$ws.Range("N12").Value2 = 478.450912
$ws.Range("O12").Value2 = 0.2983285084902258
$ws.Range("P12").Value2 = 0.2983285084902258
$ws.Range("Q12").Value2 = 21.24295468673778
$ws.Range("R12").Value2 = 191.18659218064
$ws.Range("S12").Value2 = 0.002251777459021227
$ws.Range("T12").Value2 = 0.002251777459021227

$ws.Range("G13").Value2 = 0.1331983333333333
$ws.Range("H13").Value2 = 0.399595
$ws.Range("I13").Value2 = 0.007547979475434553
$ws.Range("J13").Value2 = 0.007547979475434553
$ws.Range("O13").Value2 = 0.3227862111630279
$ws.Range("P13").Value2 = 0.3227862111630279
$ws.Range("Q13").Value2 = 22.984504202905
$ws.Range("R13").Value2 = 206.860537826145
$ws.Range("S13").Value2 = 0.002436383696811818
$ws.Range("T13").Value2 = 0.002436383696811818

$ws.Range("G14").Value2 = 0.1331983333333333
$ws.Range("H14").Value2 = 0.399595
$ws.Range("I14").Value2 = 0.007547979475434553
$ws.Range("J14").Value2 = 0.007547979475434553
$ws.Range("M14").Value2 = 74.38770566666666
$ws.Range("N14").Value2 = 223.163117
$ws.Range("O14").Value2 = 0.1391489036280481
$ws.Range("P14").Value2 = 0.1391489036280482
$ws.Range("Q14").Value2 = 9.908318415290555
$ws.Range("R14").Value2 = 89.17486573761501
$ws.Range("S14").Value2 = 0.001050293068613728
$ws.Range("T14").Value2 = 0.001050293068613728

$ws.Range("G15").Value2 = 0.1331983333333333
$ws.Range("H15").Value2 = 0.399595
$ws.Range("I15").Value2 = 0.007547979475434553
$ws.Range("J15").Value2 = 0.007547979475434553
$ws.Range("M15").Value2 = 58.41461433333333
$ws.Range("N15").Value2 = 175.243843
$ws.Range("O15").Value2 = 0.1092697975759847
$ws.Range("P15").Value2 = 0.1092697975759848
$ws.Range("Q15").Value2 = 7.780729271509444
$ws.Range("R15").Value2 = 70.026563443585
$ws.Range("S15").Value2 = 0.0008247661893884211
$ws.Range("T15").Value2 = 0.0008247661893884212

$ws.Range("G16").Value2 = 0.1331983333333333
$ws.Range("H16").Value2 = 0.399595
$ws.Range("I16").Value2 = 0.007547979475434553
$ws.Range("J16").Value2 = 0.007547979475434553
$ws.Range("M16").Value2 = 69.746216
$ws.Range("N16").Value2 = 209.238648
$ws.Range("O16").Value2 = 0.1304665791427133
$ws.Range("P16").Value2 = 0.1304665791427133
$ws.Range("Q16").Value2 = 9.290079727506667
$ws.Range("R16").Value2 = 83.61071754756001
$ws.Range("S16").Value2 = 0.0009847590615993578
$ws.Range("T16").Value2 = 0.0009847590615993578
